# "Tried to implement Penalty Reward System (unfinished)"
#
# Forecast Comparison sheet: shift Week_Start_Date forward by one week and
# overwrite MyForecast (column D) with the new, lower forecast numbers.
#
# Summary sheet: refresh the metrics that are derived from the above
# (some of these no longer reconcile perfectly with the raw data -- the
# author's own commit message flags the work as unfinished/half-wired).
#
# NOTE: every changed text cell (dates, and Summary's "numbers-as-text"
# values) must stay plain text, matching the source file's inlineStr
# cells, instead of being auto-coerced to a real number/date by Excel's
# normal typed-input parsing. We do that by switching the cell to the
# "Text" number format before writing the value, then snapping the style
# back to the default "Normal" cell style so no stray formatting lingers.

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# New Week_Start_Date values (column B, rows 2-17) -- each week's date
# shifts forward to the following week's old date (row 17 extends the
# series one more week past the old final row).
$newDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

# New MyForecast values (column D, rows 2-17).
$newForecast = @(3, 3, 3, 4, 4, 3, 3, 4, 3, 3, 3, 3, 3, 3, 3, 3)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    Set-TextValue $wsForecast.Cells.Item($row, 2) $newDates[$i]
    $wsForecast.Cells.Item($row, 4).Value = $newForecast[$i]
}

# Summary sheet updates (column B) -- kept as text, same as the source.
Set-TextValue $wsSummary.Range("B2")  "2023-01-01 to 2025-01-05"
Set-TextValue $wsSummary.Range("B4")  "67"
Set-TextValue $wsSummary.Range("B5")  "23"
Set-TextValue $wsSummary.Range("B8")  "2334 units"
Set-TextValue $wsSummary.Range("B9")  "50"
Set-TextValue $wsSummary.Range("B10") "26"
Set-TextValue $wsSummary.Range("B11") "12"
Set-TextValue $wsSummary.Range("B12") "4"
Set-TextValue $wsSummary.Range("B14") "3"
Set-TextValue $wsSummary.Range("B15") "2025-01-12"
